$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 5083
$ws.Range('I3').Value = 5297
$ws.Range('I5').Value = 5135
$ws.Range('H6').Value = 1674
$ws.Range('I6').Value = 1215
$ws.Range('I7').Value = 492
$ws.Range('I8').Value = 11085
$ws.Range('I9').Value = 5788
$ws.Range('I10').Value = 37084
$ws.Range('H11').Value = 84500
$ws.Range('I11').Value = 71471

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('I3').Value = 179
$ws.Range('I5').Value = 177
$ws.Range('I8').Value = 293
$ws.Range('I10').Value = 709
$ws.Range('I11').Value = 1765

$ws = $wb.Worksheets.Item('Museum Campus')
$ws.Range('I8').Value = 23
$ws.Range('I9').Value = 34

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('I9').Value = 71
$ws.Range('I10').Value = 624
$ws.Range('I11').Value = 1048

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range('I8').Value = 71
$ws.Range('I11').Value = 465

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('I2').Value = 38
$ws.Range('I8').Value = 62
$ws.Range('I10').Value = 256
$ws.Range('I11').Value = 483

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('I3').Value = 186
$ws.Range('I8').Value = 426
$ws.Range('I10').Value = 575
$ws.Range('I11').Value = 1695

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('I2').Value = 95
$ws.Range('I6').Value = 25
$ws.Range('I11').Value = 984

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I2').Value = 162
$ws.Range('I3').Value = 253
$ws.Range('I7').Value = 19
$ws.Range('I8').Value = 228
$ws.Range('I10').Value = 643
$ws.Range('I11').Value = 1656

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('I3').Value = 48
$ws.Range('I9').Value = 39
$ws.Range('I10').Value = 190
$ws.Range('I11').Value = 460

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I2').Value = 557
$ws.Range('I3').Value = 102
$ws.Range('I4').Value = 300
$ws.Range('I5').Value = 170
$ws.Range('I6').Value = 469
$ws.Range('I7').Value = 1765
$ws.Range('I8').Value = 3000
$ws.Range('I10').Value = 609
$ws.Range('I11').Value = 1056
$ws.Range('I13').Value = 206
$ws.Range('I14').Value = 483
$ws.Range('I15').Value = 619
$ws.Range('I16').Value = 399
$ws.Range('I18').Value = 436
$ws.Range('I19').Value = 1667
$ws.Range('I20').Value = 1399
$ws.Range('I22').Value = 269
$ws.Range('I23').Value = 897
$ws.Range('I24').Value = 295
$ws.Range('I25').Value = 272
$ws.Range('I26').Value = 173
$ws.Range('I27').Value = 805
$ws.Range('I28').Value = 61
$ws.Range('I29').Value = 2635
$ws.Range('I33').Value = 1928
$ws.Range('I34').Value = 530
$ws.Range('I36').Value = 941
$ws.Range('I37').Value = 1695
$ws.Range('I38').Value = 139
$ws.Range('I42').Value = 1721
$ws.Range('I43').Value = 717
$ws.Range('I44').Value = 714
$ws.Range('I45').Value = 138
$ws.Range('I47').Value = 572
$ws.Range('I48').Value = 1828
$ws.Range('I49').Value = 1119
$ws.Range('I50').Value = 597
$ws.Range('I51').Value = 854
$ws.Range('I52').Value = 1018
$ws.Range('I53').Value = 996
$ws.Range('I54').Value = 2437
$ws.Range('I55').Value = 834
$ws.Range('I56').Value = 397
$ws.Range('I60').Value = 454
$ws.Range('I62').Value = 34
$ws.Range('H63').Value = 2068
$ws.Range('I63').Value = 1112
$ws.Range('I64').Value = 856
$ws.Range('I65').Value = 992
$ws.Range('I66').Value = 434
$ws.Range('I67').Value = 1656
$ws.Range('I68').Value = 291
$ws.Range('I69').Value = 322
$ws.Range('I70').Value = 465
$ws.Range('I71').Value = 216
$ws.Range('I72').Value = 477
$ws.Range('I73').Value = 763
$ws.Range('I75').Value = 222
$ws.Range('I76').Value = 2027
$ws.Range('I77').Value = 278
$ws.Range('I78').Value = 1247
$ws.Range('I79').Value = 1546
$ws.Range('I80').Value = 277
$ws.Range('I81').Value = 130
$ws.Range('I82').Value = 163
$ws.Range('I83').Value = 1126
$ws.Range('I84').Value = 460
$ws.Range('I85').Value = 2481
$ws.Range('I86').Value = 771
$ws.Range('I87').Value = 245
$ws.Range('I88').Value = 621
$ws.Range('I89').Value = 1263
$ws.Range('I90').Value = 821
$ws.Range('I91').Value = 628
$ws.Range('I92').Value = 218
$ws.Range('I93').Value = 468
$ws.Range('I94').Value = 1719
$ws.Range('I95').Value = 835
$ws.Range('I96').Value = 1048
$ws.Range('I97').Value = 1030
$ws.Range('I98').Value = 914
$ws.Range('I99').Value = 984
$ws.Range('I100').Value = 153
$ws.Range('H101').Value = 84500
$ws.Range('I101').Value = 71471

$ws = $wb.Worksheets.Item('New City')
$ws.Range('I2').Value = 132
$ws.Range('I3').Value = 118
$ws.Range('I5').Value = 69
$ws.Range('I8').Value = 114
$ws.Range('I10').Value = 404
$ws.Range('I11').Value = 992

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('I10').Value = 486
$ws.Range('I11').Value = 805

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('I2').Value = 133
$ws.Range('I8').Value = 214
$ws.Range('I11').Value = 1126

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('I3').Value = 306
$ws.Range('I6').Value = 36
$ws.Range('I8').Value = 317
$ws.Range('I9').Value = 254
$ws.Range('I10').Value = 658
$ws.Range('I11').Value = 1928

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('I3').Value = 108
$ws.Range('I8').Value = 154
$ws.Range('I11').Value = 835

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('I9').Value = 14
$ws.Range('I11').Value = 222

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('I8').Value = 323
$ws.Range('I10').Value = 561
$ws.Range('I11').Value = 1546

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('I10').Value = 88
$ws.Range('I11').Value = 216

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('I2').Value = 32
$ws.Range('I5').Value = 41
$ws.Range('I8').Value = 77
$ws.Range('I10').Value = 242
$ws.Range('I11').Value = 454

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('I3').Value = 48
$ws.Range('I8').Value = 116
$ws.Range('I10').Value = 842
$ws.Range('I11').Value = 1263

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('I9').Value = 54
$ws.Range('I10').Value = 434
$ws.Range('I11').Value = 856

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('I5').Value = 60
$ws.Range('I9').Value = 104
$ws.Range('I10').Value = 1292
$ws.Range('I11').Value = 1719

$ws = $wb.Worksheets.Item('River North')
$ws.Range('I5').Value = 53
$ws.Range('I10').Value = 1468
$ws.Range('I11').Value = 2027

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('I8').Value = 51
$ws.Range('I10').Value = 285
$ws.Range('I11').Value = 434

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('I10').Value = 154
$ws.Range('I11').Value = 245

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range('I9').Value = 103
$ws.Range('I10').Value = 153

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('I10').Value = 122
$ws.Range('I11').Value = 272

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('I2').Value = 10
$ws.Range('I10').Value = 257
$ws.Range('I11').Value = 399

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('I8').Value = 156
$ws.Range('I10').Value = 752
$ws.Range('I11').Value = 1119

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('I2').Value = 220
$ws.Range('I3').Value = 320
$ws.Range('I8').Value = 499
$ws.Range('I10').Value = 901
$ws.Range('I11').Value = 2481

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('I8').Value = 152
$ws.Range('I10').Value = 650
$ws.Range('I11').Value = 1030

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('I3').Value = 65
$ws.Range('I10').Value = 439
$ws.Range('I11').Value = 834

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('I10').Value = 1707
$ws.Range('I11').Value = 2437

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('I2').Value = 48
$ws.Range('I3').Value = 51
$ws.Range('I10').Value = 421
$ws.Range('I11').Value = 763

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('I2').Value = 339
$ws.Range('I3').Value = 392
$ws.Range('I10').Value = 884
$ws.Range('I11').Value = 2635

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('I2').Value = 179
$ws.Range('I5').Value = 167
$ws.Range('I8').Value = 386
$ws.Range('I10').Value = 617
$ws.Range('I11').Value = 1667

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('I2').Value = 16
$ws.Range('I11').Value = 269

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('I10').Value = 420
$ws.Range('I11').Value = 714

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('I2').Value = 156
$ws.Range('I5').Value = 136
$ws.Range('I10').Value = 705
$ws.Range('I11').Value = 1721

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('I5').Value = 117
$ws.Range('I8').Value = 206
$ws.Range('I10').Value = 1258
$ws.Range('I11').Value = 1828

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range('I10').Value = 207
$ws.Range('I11').Value = 322

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('I2').Value = 49
$ws.Range('I11').Value = 469

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('I7').Value = 9
$ws.Range('I8').Value = 250
$ws.Range('I10').Value = 381
$ws.Range('I11').Value = 941

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('I9').Value = 152
$ws.Range('I10').Value = 206

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('I8').Value = 73
$ws.Range('I10').Value = 359
$ws.Range('I11').Value = 609

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('I2').Value = 20
$ws.Range('I6').Value = 53
$ws.Range('I10').Value = 535
$ws.Range('I11').Value = 771

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('I8').Value = 128
$ws.Range('I10').Value = 786
$ws.Range('I11').Value = 1247

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('I10').Value = 169
$ws.Range('I11').Value = 291

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('I8').Value = 76
$ws.Range('I10').Value = 280
$ws.Range('I11').Value = 619

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('I5').Value = 40
$ws.Range('I11').Value = 295

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('I3').Value = 144
$ws.Range('I8').Value = 139
$ws.Range('I11').Value = 1018

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('I3').Value = 61
$ws.Range('I9').Value = 50
$ws.Range('I10').Value = 480
$ws.Range('I11').Value = 897

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('I8').Value = 220
$ws.Range('I10').Value = 594
$ws.Range('I11').Value = 1399

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('I3').Value = 53
$ws.Range('I8').Value = 144
$ws.Range('I10').Value = 463
$ws.Range('I11').Value = 854

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('I8').Value = 153
$ws.Range('I11').Value = 628

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('I8').Value = 155
$ws.Range('I11').Value = 572

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('I5').Value = 88
$ws.Range('I8').Value = 206
$ws.Range('I11').Value = 821

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('I8').Value = 48
$ws.Range('I10').Value = 418
$ws.Range('I11').Value = 597

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('I10').Value = 247
$ws.Range('I11').Value = 468

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('I8').Value = 99
$ws.Range('I11').Value = 436

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('I8').Value = 150
$ws.Range('I10').Value = 547
$ws.Range('I11').Value = 1056

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('I5').Value = 23
$ws.Range('I11').Value = 278

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range('I8').Value = 12
$ws.Range('I10').Value = 362
$ws.Range('I11').Value = 397

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('I3').Value = 15
$ws.Range('I5').Value = 16
$ws.Range('I8').Value = 24
$ws.Range('I11').Value = 170

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('I10').Value = 319
$ws.Range('I11').Value = 557

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('I8').Value = 67
$ws.Range('I10').Value = 297
$ws.Range('I11').Value = 477

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('I2').Value = 29
$ws.Range('I8').Value = 139
$ws.Range('I10').Value = 394
$ws.Range('I11').Value = 717

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('I10').Value = 156
$ws.Range('I11').Value = 300

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('I10').Value = 179
$ws.Range('I11').Value = 277

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('I8').Value = 85
$ws.Range('I10').Value = 316
$ws.Range('I11').Value = 530

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('I2').Value = 338
$ws.Range('I3').Value = 306
$ws.Range('I5').Value = 184
$ws.Range('I8').Value = 573
$ws.Range('I9').Value = 344
$ws.Range('I10').Value = 1138
$ws.Range('I11').Value = 3000

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range('I9').Value = 76
$ws.Range('I10').Value = 138

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('I8').Value = 70
$ws.Range('I10').Value = 664
$ws.Range('I11').Value = 914

$ws = $wb.Worksheets.Item('East Village')
$ws.Range('I10').Value = 105
$ws.Range('I11').Value = 173

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('I10').Value = 105
$ws.Range('I11').Value = 218

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range('I7').Value = 23
$ws.Range('I10').Value = 163

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('I9').Value = 84
$ws.Range('I10').Value = 615
$ws.Range('I11').Value = 996

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('I2').Value = 50
$ws.Range('I9').Value = 47
$ws.Range('I10').Value = 317
$ws.Range('I11').Value = 621

$ws = $wb.Worksheets.Item('Andersonville')
$ws.Range('I9').Value = 67
$ws.Range('I10').Value = 102

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range('I9').Value = 83
$ws.Range('I10').Value = 130

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range('I9').Value = 115
$ws.Range('I10').Value = 139

$ws = $wb.Worksheets.Item('Edison Park')
$ws.Range('I10').Value = 29
$ws.Range('I11').Value = 61
